$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K10").Value = 13.62268518518517
$ws.Range("K11").Value = 13.62268518518517
$ws.Range("R11").Value = 0.9516859959956178
$ws.Range("S11").Value = 0.9753614246104579
$ws.Range("K20").Value = 12.93898809523811
$ws.Range("K21").Value = 12.93898809523811
$ws.Range("R21").Value = 0.9495669873270495
$ws.Range("S21").Value = 0.9731165936130245
$ws.Range("K24").Value = 13.46442495126706
$ws.Range("K25").Value = 13.46442495126706
$ws.Range("R25").Value = 0.9511946531936644
$ws.Range("S25").Value = 0.9748408768576692
$ws.Range("K26").Value = 15.74228395061728
$ws.Range("K27").Value = 15.74228395061728
$ws.Range("R27").Value = 0.9583158770332573
$ws.Range("S27").Value = 0.9823871532785465
$ws.Range("K34").Value = 19.79629629629628
$ws.Range("R34").Value = 0.9712571710566898
$ws.Range("S34").Value = 0.9961106502456767
$ws.Range("K35").Value = 19.79629629629628
$ws.Range("K40").Value = 20.68981481481483
$ws.Range("K41").Value = 20.68981481481483
$ws.Range("R41").Value = 0.9741566255742371
$ws.Range("S41").Value = 0.9991870829399434
$ws.Range("K48").Value = 21.28240740740739
$ws.Range("K49").Value = 21.28240740740739
$ws.Range("R49").Value = 0.9760891465058971
$ws.Range("S49").Value = 1.001237913506406
$ws.Range("K50").Value = 19.65277777777778
$ws.Range("K51").Value = 19.65277777777778
$ws.Range("R51").Value = 0.970793063583815
$ws.Range("S51").Value = 0.995618273140397
$ws.Range("K52").Value = 21.28240740740739
$ws.Range("R52").Value = 0.9760891465058971
$ws.Range("S52").Value = 1.001237913506406
$ws.Range("K53").Value = 21.28240740740739
$ws.Range("K56").Value = 19.65277777777778
$ws.Range("K65").Value = 13.62268518518517
$ws.Range("K66").Value = 13.62268518518517
$ws.Range("R66").Value = 0.9516859959956178
$ws.Range("S66").Value = 0.9753614246104579
$ws.Range("K69").Value = 19.79629629629628
$ws.Range("K70").Value = 19.79629629629628
$ws.Range("R70").Value = 0.9712571710566898
$ws.Range("S70").Value = 0.9961106502456767
$ws.Range("K71").Value = 13.75752314814816
$ws.Range("K72").Value = 13.75752314814816
$ws.Range("R72").Value = 0.9521050214763401
$ws.Range("S72").Value = 0.9758053708974481
$ws.Range("K73").Value = 19.30324074074072
$ws.Range("K74").Value = 19.30324074074072
$ws.Range("R74").Value = 0.9696645907267841
$ws.Range("S74").Value = 0.9944211305850406
$ws.Range("K75").Value = 13.00385802469133
$ws.Range("R75").Value = 0.9497676359185355
$ws.Range("S75").Value = 0.9733291418446532
$ws.Range("K76").Value = 13.00385802469133
$ws.Range("K78").Value = 13.62268518518517
$ws.Range("K79").Value = 13.62268518518517
$ws.Range("R79").Value = 0.9516859959956178
$ws.Range("S79").Value = 0.9753614246104579
$ws.Range("K80").Value = 1.791666666666668
$ws.Range("K81").Value = 1.791666666666668
$ws.Range("R81").Value = 0.9163022129108289
$ws.Range("S81").Value = 0.9379207786940652
$ws.Range("K88").Value = 13.62268518518517
$ws.Range("R88").Value = 0.9516859959956178
$ws.Range("S88").Value = 0.9753614246104579
$ws.Range("K89").Value = 13.62268518518517
$ws.Range("K90").Value = 5.486111111111112
$ws.Range("R90").Value = 0.9270655773901523
$ws.Range("S90").Value = 0.9492998859749143
$ws.Range("K91").Value = 5.486111111111112
$ws.Range("K92").Value = 1.925925925925943
$ws.Range("R92").Value = 0.9166889846297158
$ws.Range("S92").Value = 0.9383295263284442
$ws.Range("K93").Value = 1.925925925925943
$ws.Range("K102").Value = 13.75752314814816
$ws.Range("R102").Value = 0.9521050214763401
$ws.Range("S102").Value = 0.9758053708974481
$ws.Range("K103").Value = 13.75752314814816
$ws.Range("K104").Value = 13.62268518518517
$ws.Range("R104").Value = 0.9516859959956178
$ws.Range("S104").Value = 0.9753614246104579
$ws.Range("K105").Value = 13.62268518518517
$ws.Range("K110").Value = 1.791666666666668
$ws.Range("K111").Value = 1.791666666666668
$ws.Range("R111").Value = 0.9163022129108289
$ws.Range("S111").Value = 0.9379207786940652
$ws.Range("K112").Value = 21.19907407407406
$ws.Range("R112").Value = 0.9758169225763468
$ws.Range("S112").Value = 1.000949006909155
$ws.Range("K113").Value = 21.19907407407406

Write-Output "Applied 93 cell updates"
